$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Update cell values for rows 8-10 on all three sheets to reflect
# the reordering of cf41b5ba / cf964e1e / d7976345 entries and their
# associated status text.
# ---------------------------------------------------------------

$ws = $wb.Worksheets.Item(1)
$ws.Range("A8").Value = "cf964e1e-f792-415a-bbf2-d47895a02c50.md"
$ws.Range("B8").Value = "In Translation"
$ws.Range("C8").Value = "In Translation"
$ws.Range("A9").Value = "d7976345-19df-48d5-bcd3-479ca80dd078.md"
$ws.Range("B9").Value = "In Translation"
$ws.Range("C9").Value = "In Translation"
$ws.Range("A10").Value = "cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md"
$ws.Range("B10").Value = "Ready for handoff"
$ws.Range("C10").Value = "Ready for handoff"

$ws = $wb.Worksheets.Item(2)
$ws.Range("A8").Value = "cf964e1e-f792-415a-bbf2-d47895a02c50.md"
$ws.Range("B8").Value = "In Translation"
$ws.Range("A9").Value = "d7976345-19df-48d5-bcd3-479ca80dd078.md"
$ws.Range("B9").Value = "In Translation"
$ws.Range("A10").Value = "cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md"
$ws.Range("B10").Value = "Ready for handoff"

$ws = $wb.Worksheets.Item(3)
$ws.Range("A8").Value = "cf964e1e-f792-415a-bbf2-d47895a02c50.md"
$ws.Range("B8").Value = "In Translation"
$ws.Range("A9").Value = "d7976345-19df-48d5-bcd3-479ca80dd078.md"
$ws.Range("B9").Value = "In Translation"
$ws.Range("A10").Value = "cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md"
$ws.Range("B10").Value = "Ready for handoff"

# ---------------------------------------------------------------
# Rebuild hyperlinks on each sheet. The underlying link targets
# are unchanged; only which display text/cell shows which link
# moves, matching the new row order above. Because this runtime
# always appends hyperlinks rather than editing them in place, we
# clear each sheet's hyperlinks and recreate the full set in the
# original relationship order with the corrected display text.
# ---------------------------------------------------------------

$ws = $wb.Worksheets.Item(1)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/1364424f-b177-40c0-a6d0-00face97328a.md", "", "", "1364424f-b177-40c0-a6d0-00face97328a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ad2540ea7ff59cd19be6996f1aedaebed8bb62bb/e2e/15975468-7825-45d1-a84f-ccf1963b9399.md", "", "", "15975468-7825-45d1-a84f-ccf1963b9399.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ad2540ea7ff59cd19be6996f1aedaebed8bb62bb/e2e/863e8a40-f220-4138-bc79-3a9799e60980.md", "", "", "863e8a40-f220-4138-bc79-3a9799e60980.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/b48ee028-de37-4f31-bd6a-56789eb1ad77.md", "", "", "b48ee028-de37-4f31-bd6a-56789eb1ad77.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/eb922731-4057-4c54-b814-2a35b7a1adad.md", "", "", "eb922731-4057-4c54-b814-2a35b7a1adad.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/51355f24322ca0b50366d2d9ca66c3f5d2408931/e2e/93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.md", "", "", "93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md", "", "", "cf964e1e-f792-415a-bbf2-d47895a02c50.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/cf964e1e-f792-415a-bbf2-d47895a02c50.md", "", "", "d7976345-19df-48d5-bcd3-479ca80dd078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/d7976345-19df-48d5-bcd3-479ca80dd078.md", "", "", "cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A11"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/.localization-config", "", "", ".localization-config") | Out-Null

$ws = $wb.Worksheets.Item(2)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/1364424f-b177-40c0-a6d0-00face97328a.md", "", "", "1364424f-b177-40c0-a6d0-00face97328a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35a49c93accf41fb066b7df7b7f36c5174b90cfe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/1364424f-b177-40c0-a6d0-00face97328a.12b5ac0775be1a0162339e2a0748007934e6b375.zh-cn.xlf", "", "", "1364424f-b177-40c0-a6d0-00face97328a.12b5ac0775be1a0162339e2a0748007934e6b375.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b69ea5ea17aa98a2513fc1492525052c6382a5c7/e2e/1364424f-b177-40c0-a6d0-00face97328a.md", "", "", "1364424f-b177-40c0-a6d0-00face97328a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/64653f5dfa21e41c9e29b6b0114d9cb234300050/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/mt/1364424f-b177-40c0-a6d0-00face97328a.12b5ac0775be1a0162339e2a0748007934e6b375.zh-cn.xlf", "", "", "1364424f-b177-40c0-a6d0-00face97328a.12b5ac0775be1a0162339e2a0748007934e6b375.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ad2540ea7ff59cd19be6996f1aedaebed8bb62bb/e2e/15975468-7825-45d1-a84f-ccf1963b9399.md", "", "", "15975468-7825-45d1-a84f-ccf1963b9399.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e72054e06fd4434429110b735690cce668cb2f43/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/15975468-7825-45d1-a84f-ccf1963b9399.7990d346f556bbfd654f00a11a20d84a33e6ebed.zh-cn.xlf", "", "", "15975468-7825-45d1-a84f-ccf1963b9399.7990d346f556bbfd654f00a11a20d84a33e6ebed.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1804b2d034cf67a78cb973e94a2d3874d36c493b/e2e/15975468-7825-45d1-a84f-ccf1963b9399.md", "", "", "15975468-7825-45d1-a84f-ccf1963b9399.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3eb58d2e5e5ccad83199b6eb6064e633f3b4db3f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/15975468-7825-45d1-a84f-ccf1963b9399.7990d346f556bbfd654f00a11a20d84a33e6ebed.zh-cn.xlf", "", "", "15975468-7825-45d1-a84f-ccf1963b9399.7990d346f556bbfd654f00a11a20d84a33e6ebed.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ad2540ea7ff59cd19be6996f1aedaebed8bb62bb/e2e/863e8a40-f220-4138-bc79-3a9799e60980.md", "", "", "863e8a40-f220-4138-bc79-3a9799e60980.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e72054e06fd4434429110b735690cce668cb2f43/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/863e8a40-f220-4138-bc79-3a9799e60980.28497e9433b5d758641c5a651c1062562c1864c5.zh-cn.xlf", "", "", "863e8a40-f220-4138-bc79-3a9799e60980.28497e9433b5d758641c5a651c1062562c1864c5.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1804b2d034cf67a78cb973e94a2d3874d36c493b/e2e/863e8a40-f220-4138-bc79-3a9799e60980.md", "", "", "863e8a40-f220-4138-bc79-3a9799e60980.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3eb58d2e5e5ccad83199b6eb6064e633f3b4db3f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/863e8a40-f220-4138-bc79-3a9799e60980.28497e9433b5d758641c5a651c1062562c1864c5.zh-cn.xlf", "", "", "863e8a40-f220-4138-bc79-3a9799e60980.28497e9433b5d758641c5a651c1062562c1864c5.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/b48ee028-de37-4f31-bd6a-56789eb1ad77.md", "", "", "b48ee028-de37-4f31-bd6a-56789eb1ad77.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35a49c93accf41fb066b7df7b7f36c5174b90cfe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/b48ee028-de37-4f31-bd6a-56789eb1ad77.7118de5c76beb28c9bfab1e282a576c6468937bc.zh-cn.xlf", "", "", "b48ee028-de37-4f31-bd6a-56789eb1ad77.7118de5c76beb28c9bfab1e282a576c6468937bc.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b69ea5ea17aa98a2513fc1492525052c6382a5c7/e2e/b48ee028-de37-4f31-bd6a-56789eb1ad77.md", "", "", "b48ee028-de37-4f31-bd6a-56789eb1ad77.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/64653f5dfa21e41c9e29b6b0114d9cb234300050/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/mt/b48ee028-de37-4f31-bd6a-56789eb1ad77.7118de5c76beb28c9bfab1e282a576c6468937bc.zh-cn.xlf", "", "", "b48ee028-de37-4f31-bd6a-56789eb1ad77.7118de5c76beb28c9bfab1e282a576c6468937bc.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/eb922731-4057-4c54-b814-2a35b7a1adad.md", "", "", "eb922731-4057-4c54-b814-2a35b7a1adad.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35a49c93accf41fb066b7df7b7f36c5174b90cfe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/eb922731-4057-4c54-b814-2a35b7a1adad.f4fcb4a445ae7496b96f83bd6cf26c1c131018e4.zh-cn.xlf", "", "", "eb922731-4057-4c54-b814-2a35b7a1adad.f4fcb4a445ae7496b96f83bd6cf26c1c131018e4.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/51355f24322ca0b50366d2d9ca66c3f5d2408931/e2e/93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.md", "", "", "93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5263e1af4bbb1034f022ef4bd1e2d29a5039bb88/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.76e0d82dd90a63c80294ab815928575f4d4dff3b.zh-cn.xlf", "", "", "93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.76e0d82dd90a63c80294ab815928575f4d4dff3b.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md", "", "", "cf964e1e-f792-415a-bbf2-d47895a02c50.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35a49c93accf41fb066b7df7b7f36c5174b90cfe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.9c8087550961789a7eba406c4c205f2f28e4bf3e.zh-cn.xlf", "", "", "cf964e1e-f792-415a-bbf2-d47895a02c50.e9b5800a0396377cdf80d21a46e369c3f770dee3.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/cf964e1e-f792-415a-bbf2-d47895a02c50.md", "", "", "d7976345-19df-48d5-bcd3-479ca80dd078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35a49c93accf41fb066b7df7b7f36c5174b90cfe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/cf964e1e-f792-415a-bbf2-d47895a02c50.e9b5800a0396377cdf80d21a46e369c3f770dee3.zh-cn.xlf", "", "", "d7976345-19df-48d5-bcd3-479ca80dd078.c3fc8a37fc262f76e596b699ade11668d231fbb4.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/d7976345-19df-48d5-bcd3-479ca80dd078.md", "", "", "cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35a49c93accf41fb066b7df7b7f36c5174b90cfe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/d7976345-19df-48d5-bcd3-479ca80dd078.c3fc8a37fc262f76e596b699ade11668d231fbb4.zh-cn.xlf", "", "", "cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.9c8087550961789a7eba406c4c205f2f28e4bf3e.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A11"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/.localization-config", "", "", ".localization-config") | Out-Null

$ws = $wb.Worksheets.Item(3)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/1364424f-b177-40c0-a6d0-00face97328a.md", "", "", "1364424f-b177-40c0-a6d0-00face97328a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f6133b26f10644e717fb4dcd66f711f74bf8626/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/1364424f-b177-40c0-a6d0-00face97328a.12b5ac0775be1a0162339e2a0748007934e6b375.de-de.xlf", "", "", "1364424f-b177-40c0-a6d0-00face97328a.12b5ac0775be1a0162339e2a0748007934e6b375.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e1c9ae631e306999edca1797b39354ad1e91ad9b/e2e/1364424f-b177-40c0-a6d0-00face97328a.md", "", "", "1364424f-b177-40c0-a6d0-00face97328a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a4c3f93ee08d2da800f89a8b558779c0dd3f9eba/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/mt/1364424f-b177-40c0-a6d0-00face97328a.12b5ac0775be1a0162339e2a0748007934e6b375.de-de.xlf", "", "", "1364424f-b177-40c0-a6d0-00face97328a.12b5ac0775be1a0162339e2a0748007934e6b375.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ad2540ea7ff59cd19be6996f1aedaebed8bb62bb/e2e/15975468-7825-45d1-a84f-ccf1963b9399.md", "", "", "15975468-7825-45d1-a84f-ccf1963b9399.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/004b0a76b510d3225ef12377eb26f79f0abb554d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/15975468-7825-45d1-a84f-ccf1963b9399.7990d346f556bbfd654f00a11a20d84a33e6ebed.de-de.xlf", "", "", "15975468-7825-45d1-a84f-ccf1963b9399.7990d346f556bbfd654f00a11a20d84a33e6ebed.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/41dc0c81a7a89a3ba6e019a554a99df0a4e26879/e2e/15975468-7825-45d1-a84f-ccf1963b9399.md", "", "", "15975468-7825-45d1-a84f-ccf1963b9399.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/685fb9c9f04243f1e74053f2015653b54026aa71/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/15975468-7825-45d1-a84f-ccf1963b9399.7990d346f556bbfd654f00a11a20d84a33e6ebed.de-de.xlf", "", "", "15975468-7825-45d1-a84f-ccf1963b9399.7990d346f556bbfd654f00a11a20d84a33e6ebed.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ad2540ea7ff59cd19be6996f1aedaebed8bb62bb/e2e/863e8a40-f220-4138-bc79-3a9799e60980.md", "", "", "863e8a40-f220-4138-bc79-3a9799e60980.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/004b0a76b510d3225ef12377eb26f79f0abb554d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/863e8a40-f220-4138-bc79-3a9799e60980.28497e9433b5d758641c5a651c1062562c1864c5.de-de.xlf", "", "", "863e8a40-f220-4138-bc79-3a9799e60980.28497e9433b5d758641c5a651c1062562c1864c5.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/41dc0c81a7a89a3ba6e019a554a99df0a4e26879/e2e/863e8a40-f220-4138-bc79-3a9799e60980.md", "", "", "863e8a40-f220-4138-bc79-3a9799e60980.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/685fb9c9f04243f1e74053f2015653b54026aa71/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/863e8a40-f220-4138-bc79-3a9799e60980.28497e9433b5d758641c5a651c1062562c1864c5.de-de.xlf", "", "", "863e8a40-f220-4138-bc79-3a9799e60980.28497e9433b5d758641c5a651c1062562c1864c5.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/b48ee028-de37-4f31-bd6a-56789eb1ad77.md", "", "", "b48ee028-de37-4f31-bd6a-56789eb1ad77.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f6133b26f10644e717fb4dcd66f711f74bf8626/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/b48ee028-de37-4f31-bd6a-56789eb1ad77.7118de5c76beb28c9bfab1e282a576c6468937bc.de-de.xlf", "", "", "b48ee028-de37-4f31-bd6a-56789eb1ad77.7118de5c76beb28c9bfab1e282a576c6468937bc.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e1c9ae631e306999edca1797b39354ad1e91ad9b/e2e/b48ee028-de37-4f31-bd6a-56789eb1ad77.md", "", "", "b48ee028-de37-4f31-bd6a-56789eb1ad77.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a4c3f93ee08d2da800f89a8b558779c0dd3f9eba/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/mt/b48ee028-de37-4f31-bd6a-56789eb1ad77.7118de5c76beb28c9bfab1e282a576c6468937bc.de-de.xlf", "", "", "b48ee028-de37-4f31-bd6a-56789eb1ad77.7118de5c76beb28c9bfab1e282a576c6468937bc.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/eb922731-4057-4c54-b814-2a35b7a1adad.md", "", "", "eb922731-4057-4c54-b814-2a35b7a1adad.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f6133b26f10644e717fb4dcd66f711f74bf8626/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/eb922731-4057-4c54-b814-2a35b7a1adad.f4fcb4a445ae7496b96f83bd6cf26c1c131018e4.de-de.xlf", "", "", "eb922731-4057-4c54-b814-2a35b7a1adad.f4fcb4a445ae7496b96f83bd6cf26c1c131018e4.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/51355f24322ca0b50366d2d9ca66c3f5d2408931/e2e/93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.md", "", "", "93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e19dea18082fb8e1dc50edba8de5b0c0abba88b7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.76e0d82dd90a63c80294ab815928575f4d4dff3b.de-de.xlf", "", "", "93b8a23c-5e8f-4f3d-b6f4-ff2e67a83e08.76e0d82dd90a63c80294ab815928575f4d4dff3b.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md", "", "", "cf964e1e-f792-415a-bbf2-d47895a02c50.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f6133b26f10644e717fb4dcd66f711f74bf8626/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.9c8087550961789a7eba406c4c205f2f28e4bf3e.de-de.xlf", "", "", "cf964e1e-f792-415a-bbf2-d47895a02c50.e9b5800a0396377cdf80d21a46e369c3f770dee3.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/cf964e1e-f792-415a-bbf2-d47895a02c50.md", "", "", "d7976345-19df-48d5-bcd3-479ca80dd078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f6133b26f10644e717fb4dcd66f711f74bf8626/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/cf964e1e-f792-415a-bbf2-d47895a02c50.e9b5800a0396377cdf80d21a46e369c3f770dee3.de-de.xlf", "", "", "d7976345-19df-48d5-bcd3-479ca80dd078.c3fc8a37fc262f76e596b699ade11668d231fbb4.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/e2e/d7976345-19df-48d5-bcd3-479ca80dd078.md", "", "", "cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f6133b26f10644e717fb4dcd66f711f74bf8626/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/d7976345-19df-48d5-bcd3-479ca80dd078.c3fc8a37fc262f76e596b699ade11668d231fbb4.de-de.xlf", "", "", "cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.9c8087550961789a7eba406c4c205f2f28e4bf3e.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A11"), "https://github.com/OpenLocalizationTest/oltest/blob/ae0bf0c5b7a61c3838fd72aade222e46452f0532/.localization-config", "", "", ".localization-config") | Out-Null

Write-Host "Report regenerated for archive."
